$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Policies008")
$ws1.Range("A2").Value = "PLT008a"
$ws1.Range("B2").Value = "8a"

$ws2 = $wb.Worksheets.Item("Policies008_1")
$ws2.Range("A2").Value = "PLT008b"
$ws2.Range("B2").Value = "8b"

$ws3 = $wb.Worksheets.Item("Policies008_3")
$ws3.Range("A2").Value = "PLT008d"
$ws3.Range("B2").Value = "8d"

$ws4 = $wb.Worksheets.Item("Policies008_4")
$ws4.Range("A2").Value = "PLT008E"
$ws4.Range("B2").Value = "8E"
